$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'اجسام مالوری از چه چیزی مشتق شده‌اند؟'
$ws.Range("A3").Value = 'یک بیمار پسر 13 ساله از درد در ناحیه دندان های جلویی فک بالا شکایت دارد و سابقه زمین خوردن در حین بازی را گزارش می‌دهد. معاینه داخل دهانی نشان‌دهنده درجه i لق شدن دندان پیشین مرکزی چپ فک بالا است. یافته‌های رادیولوژیک نشان‌دهنده شکستگی ریشه دندان 21 است. کدام یک از گزینه‌های زیر می‌تواند کم‌ترین شکل مطلوب ترمیم برای این دندان باشد؟'
$ws.Range("A4").Value = 'کدام یک از داروهای زیر برای قطع و پیشگیری از تاکیکاردی فوق بطنی پاروکسیسمال استفاده می‌شود؟'
$ws.Range("A5").Value = 'شایع‌ترین علت مرگ در کودک مبتلا به دیفتری چیست؟'
$ws.Range("A6").Value = 'فِنسیکلیدین با همه‌ی موارد زیر شناخته می‌شود به جز:'
$ws.Range("A7").Value = 'یک جوان توسط پلیس از ایستگاه راه‌آهن آورده شده است. او نامربوط صحبت می‌کند. خشکی دهان، پوست گرم، مردمک‌های گشاد، راه رفتن نامتعادل و گفتار نامفهوم دارد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A8").Value = 'عصب دهی عضله گلوتئوس مدیوس توسط کدام عصب انجام می‌شود؟'
$ws.Range("A9").Value = 'کدام یک از داروهای زیر باعث ایجاد تغییرات در شبکیه نمی‌شود؟'
$ws.Range("A10").Value = 'پدیده رینود برود در کدام یک دیده می‌شود؟'
$ws.Range("A11").Value = 'کدام یک از عبارات زیر در مورد ضد انعقادهای خوراکی نادرست است؟'
$ws.Range("A12").Value = 'داروهای ضد روان‌پریشی از طریق چه مکانیسمی عمل می‌کنند؟ سپتامبر 2012'
$ws.Range("A13").Value = 'حلقه عمیق در کجا قرار دارد؟'
$ws.Range("A14").Value = 'یک بیمار مرد 20 ساله با درد در هنگام حرکت مراجعه می‌کند. تصویربرداری اشعه ایکس از مفصل زانو، یک ضایعه لیتیک در انتهای فوقانی درشت‌نی را نشان می‌دهد. تشخیص احتمالی چیست؟'
$ws.Range("A15").Value = 'کدام یک از عبارات زیر در مورد انسفالیت ژاپنی صحیح نیست؟'
$ws.Range("A16").Value = 'ظرفیت طبیعی لگنچه کلیه چقدر است؟'
$ws.Range("A17").Value = 'ویژگی های سندرم سروتونین مرتبط با ssri و maoi ها کدامند؟'
$ws.Range("A18").Value = 'ویژگی‌های سندرم پاسخ التهابی سیستمیک (sirs) شامل موارد زیر است به جز-'
$ws.Range("A19").Value = 'کدام یک از عوامل زیر با پیش‌آگهی ضعیف در لوسمی لنفوبلاستیک حاد (all) در کودکان مرتبط است؟'
$ws.Range("A20").Value = 'اجسام پساموما در همه موارد زیر دیده می‌شوند به جز؟'
$ws.Range("A21").Value = 'یک پسر 7 ساله بیمار است و تب دارد، نیاز کالری او چقدر است؟'
$ws.Range("A22").Value = 'کدام یک از موارد زیر از درمان‌های جسمی مورد استفاده در افسردگی نیست؟'
$ws.Range("A23").Value = 'کدام یک از هورمون‌های زیر مشتق از اسید آمینه نیست؟'
$ws.Range("A24").Value = 'ال-دوپا در درمان پارکینسون با کاربیدوپا ترکیب می‌شود زیرا -'
$ws.Range("A25").Value = 'سن معمول بروز تشنج ناشی از تب در کودکان چیست؟'
$ws.Range("A26").Value = 'تمام عضلات تنار دست توسط عصب مدیان عصب‌دهی می‌شوند به جز'
$ws.Range("A27").Value = 'در یک بیمار که اسپیرومتری انجام داده است، نسبت fev1/fvc کمتر از 0.7 و fev1 بین 30% تا 50% مقدار پیش‌بینی شده است. تشخیص این بیمار بر اساس معیارهای gold چیست؟'
$ws.Range("A28").Value = 'یک مرد 78 ساله با زردی همراه با خارش افزایش یافته و کیسه صفرای قابل لمس مراجعه می‌کند. تشخیص احتمالی چیست؟'
$ws.Range("A29").Value = 'همه موارد زیر در مورد پیوند پوست صحیح هستند، به جز:'
$ws.Range("A30").Value = 'در تروما از نوع "jiggling"، نیروی اکلوژن باعث چه چیزی می‌شود؟'
$ws.Range("A31").Value = 'ترومبوز سینوس کاورنوس پس از سینوزیت منجر به همه علائم زیر به جز کدام یک می‌شود؟'
$ws.Range("A32").Value = 'آنتاگونیست بتا 1 مورد استفاده در نارسایی احتقانی قلب کدام است؟'
$ws.Range("A33").Value = 'در روان‌پریشی کورساکوف همه موارد زیر دیده می‌شوند به جز -  
الف) از دست دادن حافظه دور  
ب) از دست دادن عملکرد فکری اما حفظ حافظه  
ج) عدم بینش، ناتوانی در درک ناتوانی  
د) حالت برگشت‌پذیر'
$ws.Range("A34").Value = 'عامل ایجاد کننده اپیدمی دراپسی چیست؟'
$ws.Range("A35").Value = 'یک کودک ۷ ماهه با تأخیر در رشد عصبی، از دست دادن مهارت‌های حرکتی، افزایش واکنش وحشت به صدا، مشکلات تغذیه‌ای و عدم رشد مناسب مراجعه کرده است. در معاینه، ماکروسفالی همراه با میوکلونوس و اسپاستیسیتی مشاهده شد. هیچ شواهدی از کدورت قرنیه یا آنژیوکراتوم وجود نداشت. در بررسی‌های بیشتر، سی‌تی اسکن هیچ نشانه‌ای از هپاتوسپلنومگالی یا درگیری قلبی نشان نداد. معاینه چشم‌پزشکی نشان‌دهنده‌ی کدام کمبود آنزیمی مسئول این علائم است؟'
$ws.Range("A36").Value = 'کارسینوم سلول سباسه پلک ها از کجا منشأ می گیرد؟'
$ws.Range("A37").Value = 'در کبد، پیش‌سازهای اجسام کتونی کدامند؟'
$ws.Range("A38").Value = 'ناحیه اینگوئینال نشان داده شده در زیر در کدام بیماری دیده می‌شود؟'
$ws.Range("A39").Value = 'مرحله کبدی اریتروپویز در چه زمانی رخ می‌دهد؟'
$ws.Range("A40").Value = 'مزایای بستن دیرهنگام بند ناف شامل همه موارد زیر است به جز؟'
$ws.Range("A41").Value = 'یک مرد ۲۵ ساله با تب و بی‌حالی از ۲ هفته پیش، آرتریت مفصل مچ پا و ندول‌های اریتماتوز حساس روی ساق پا. تشخیص چیست؟'
$ws.Range("A42").Value = 'همه موارد زیر در مورد شاهد خبره صحیح است به جز'
$ws.Range("A43").Value = 'تومور کیستیک سوپراتنتوریال کدام است؟'
$ws.Range("A44").Value = 'ریشه های عصبی رادیال:'
$ws.Range("A45").Value = 'ظاهر پرتوهای خورشیدی مشخصه کدام است؟'
$ws.Range("A46").Value = 'نسبت دور کمر به باسن که نشان دهنده چاقی در مردان است کدام است؟'
$ws.Range("A47").Value = 'ضد سرطان‌ها کدامند؟'
$ws.Range("A48").Value = 'کارسینوم سلول سنگفرشی مهاجم چگونه از کارسینوم درجا متمایز می‌شود؟'
$ws.Range("A49").Value = 'محل زخم پارگی که شبیه زخم بریده شده است کدام است؟'
$ws.Range("A50").Value = 'ابزار مورد استفاده در فیکساسیون داخلی استخوان های کوچک و همچنین در تکنیک ایلیزاروف؟'
$ws.Range("A51").Value = 'عصب پتروسال بزرگ از کدام قسمت تشکیل شده است؟'
$ws.Range("A52").Value = 'اولین علامت بیماری تی ساکس چیست؟'
$ws.Range("A53").Value = 'مثلث برایانت در تشخیص موارد زیر مفید است به جز:'
$ws.Range("A54").Value = 'کدام یک از موارد زیر یک ارگانیسم اجباری درون سلولی است؟'
$ws.Range("A55").Value = 'در درمان ارتودنسی، برای جلوگیری از آسیب به بافت‌ها، نیروهای اعمال شده معمولاً نباید از چه حدی تجاوز کنند؟'
$ws.Range("A56").Value = 'جذب گلوکز از روده به سلول‌های مخاطی روده با na+,k+ atpase همراه است. در مقابل، حرکت گلوکز از سلول‌های اپیتلیال روده به جریان خون زیرمخاطی از طریق انتقال غیرفعال انجام می‌شود. با توجه به این حقایق، کدام یک از عبارات زیر می‌تواند در زمان‌های مختلف درست باشد؟'
$ws.Range("A57").Value = 'تمام موارد زیر به جز یکی از عوارض لنفادم مزمن اندام هستند -'
$ws.Range("A58").Value = 'جداسازی قطعات dna توسط کدام روش انجام می‌شود؟'
$ws.Range("A59").Value = 'شایع‌ترین شکل هیستری تجزیه‌ای کدام است؟'
$ws.Range("A60").Value = 'بیمار تحت عکسبرداری قفسه سینه پس از یک تصادف خودرویی، با مدیاستینوم بزرگ شده و آدنوپاتی هیلار دو طرفه و پاراتراکئال راست مواجه شده است. بیمار بدون علامت بوده است، اما معاینه دقیق نشان دهنده غدد لنفاوی گردنی بزرگ شده است. این غده بیوپسی شده و درگیری توسط گرانولومای کوچک و غیر کازئوز نشان داده می‌شود. درون گرانولوماها، گاهی سلول‌های غول‌آسا با اجسام ستاره‌ای دیده می‌شوند. این‌ها به احتمال زیاد کدام یک از موارد زیر هستند؟'
$ws.Range("A61").Value = 'شفافیت رادیولوژیک مشاهده شده بین دندان های پیشین مرکزی فک بالا معمولاً به دلیل کدام مورد است؟'
$ws.Range("A62").Value = 'در طول گرسنگی، عضله از چه چیزی استفاده می‌کند؟'
$ws.Range("A63").Value = 'نقطه مک برنی کجاست؟ مارس 2013 (a, c, d, e)'
$ws.Range("A64").Value = 'شستشوی معده در کدام مورد منع مصرف دارد؟'
$ws.Range("A65").Value = 'یک مرد 30 ساله در یک تصادف موتورسیکلت که اندام تحتانی راست او زیر موتورسیکلت گیر کرده بود، دچار آسیب عصب گلوتئال فوقانی شد. او در بخش اورژانس تثبیت شده است. بعداً معاینه می‌شود و الگوی راه رفتن اردکی و علامت ترندلنبرگ مثبت نشان می‌دهد. کدام یک از یافته‌های فیزیکی زیر به احتمال زیاد در این بیمار مشاهده می‌شود؟'
$ws.Range("A66").Value = 'یک زن ۳۰ ساله با کاهش شنوایی هدایتی پیشرونده دوطرفه مراجعه می‌کند. شایع‌ترین علت این حالت چیست؟'
$ws.Range("A67").Value = 'ivp باید با احتیاط در چه موردی انجام شود؟'
$ws.Range("A68").Value = 'کدام داروی ضد سل باعث نوریت اپتیک می‌شود؟'
$ws.Range("A69").Value = 'میزبان واسط انسان برای کدام یک از موارد زیر است؟'
$ws.Range("A70").Value = 'واکنش مردمک به نور در آسیب به کدام یک همچنان طبیعی است؟'
$ws.Range("A71").Value = 'عبارت صحیح در میان گزینه‌های زیر کدام است؟'
$ws.Range("A72").Value = 'پروفیل انعقادی در یک دختر ۱۳ ساله با منوراژی مبتلا به بیماری فون ویلبراند چیست؟'
$ws.Range("A73").Value = 'کدام یک از موارد زیر در مورد اختلال نشخوار صحیح نیست؟'
$ws.Range("A74").Value = 'علت نکروز کازئوز در سل چیست؟'
$ws.Range("A75").Value = 'همه موارد زیر در مورد مردمک تونیک آدی صحیح است به جز:'
$ws.Range("A76").Value = 'یک خانم با 7 هفته آمنوره و لکه‌بینی واژینال خفیف مراجعه کرده است. در سونوگرافی ترانس واژینال، طول جنین (crl) 5 میلی‌متر و ساک حاملگی با شکل‌گیری مناسب مشاهده شد که سن حاملگی محاسبه‌شده 5.6 هفته بود. اقدام بعدی چیست؟'
$ws.Range("A77").Value = 'یک مرد 60 ساله با تغییر رنگ، ضخیم شدن و حفره‌دار شدن 2 ناخن دست و یک ناخن پا مراجعه کرده است. کدام یک از موارد زیر می‌تواند در سریع‌ترین زمان تشخیص را قطعی کند؟'
$ws.Range("A78").Value = 'در طی عمل تیروئید، عصبی که همراه با شریان تیروئید فوقانی حرکت می‌کند آسیب می‌بیند. چه عواقبی ممکن است رخ دهد؟'
$ws.Range("A79").Value = 'دیورز ناشی از الکل به دلیل کدام مورد است؟'
$ws.Range("A80").Value = 'سوماتواستاتین توسط کدام سلول تولید می‌شود؟'
$ws.Range("A81").Value = 'فرورفتگی کارهارت در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range("A82").Value = 'کاتاراکت کاپولیفرم از کجا شروع می‌شود؟'
$ws.Range("A83").Value = 'تکنیک رنگ‌آمیزی استفاده شده در اینجا چیست؟'
$ws.Range("A84").Value = 'در مورد کارسینوم فیبرولاملار کبد همه موارد زیر صحیح است به جز-'
$ws.Range("A85").Value = 'کدام یک از داروهای زیر در درمان عفونت ناشی از مایکوباکتریوم اویوم کمپلکس مفید است؟'
$ws.Range("A86").Value = 'در زردی غیر همولیتیک، اوروبیلینوژن در کدام مورد دیده می‌شود؟'
$ws.Range("A87").Value = 'در تنگی شدید آئورت، یافته صحیح کدام است؟'
$ws.Range("A88").Value = 'مکانیسم حمل و نقلی که به پروتئین حامل نیاز دارد اما به انرژی نیاز ندارد -'
$ws.Range("A89").Value = 'هورمون پاراتیروئید با اثر بر کدام بخش، جذب کلسیم را افزایش می‌دهد؟'
$ws.Range("A90").Value = 'یک داروی اسیدی دارای مقدار pka برابر با 4.5 است. درصد داروی یونیزه شده در ph 3.5 چقدر است؟'
$ws.Range("A91").Value = 'نارسایی حاد و مزمن کلیه را می‌توان با کدام مورد از هم تفکیک کرد؟'
$ws.Range("A92").Value = 'یک مرد 23 ساله با تب، کاهش وزن، بی‌حالی، درد شکم و میالژی مراجعه می‌کند. بررسی‌ها نشان می‌دهد که بیمار به پلی‌آرتریت ندوزا مبتلا است. کدام یک از موارد زیر با این نوع واسکولیت مرتبط است؟'
$ws.Range("A93").Value = 'کدام یک از تومورهای تخمدانی زیر بیشترین حساسیت را به پرتودرمانی دارد؟'
$ws.Range("A94").Value = 'ایمنی غیرفعال مادرزادی در کدام یک از موارد زیر یافت نمی‌شود؟'
$ws.Range("A95").Value = 'فلاپ میوکوتانئوس شامل کدام بافت‌ها می‌شود؟'
$ws.Range("A96").Value = 'طول عصب بینایی چقدر است؟'
$ws.Range("A97").Value = 'در چرخش شدید پا به داخل، شکستگی توبروزیته متاتارس پنجم به دلیل کشش بیش از حد کدام تاندون رخ می‌دهد؟'
$ws.Range("A98").Value = 'پس از کشیدن دندانی که در ناحیه اطراف ریشه رادیولوسنسی نشان می‌دهد'
$ws.Range("A99").Value = 'کدام سینوس بر روی فالکس مخچه ای قرار دارد؟'
$ws.Range("A100").Value = 'عضله استاپدیوس که استخوان رکابی را حرکت می‌دهد، توسط کدام عصب عصب‌دهی می‌شود؟'
$ws.Range("A101").Value = 'شایع‌ترین گروه سنی برای بوتولیسم نوزادی کدام است؟'
$ws.Range("A102").Value = 'در مورد انسولینوما کدام گزینه صحیح است؟'
$ws.Range("A103").Value = 'کدام استخوان بیشتر تحت تأثیر پوکی استخوان ناشی از گلوکوکورتیکوئید قرار می‌گیرد؟'
$ws.Range("A104").Value = 'تیوپنتون برای بیهوشی القایی استفاده می‌شود، زیرا -'
$ws.Range("A105").Value = 'رتینوپاتی پوشر در بیماران مبتلا به کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A106").Value = 'کدام یک از علائم بالینی زیر در بیماری چشمی گریوز دیده نمی‌شود؟'
$ws.Range("A107").Value = 'حساسترین آزمایش برای سندرم تونل کارپال'
$ws.Range("A108").Value = 'همه موارد زیر می‌توانند به ایجاد گلوکوم با زاویه بسته predispose کنند، به جز:'
$ws.Range("A109").Value = 'داروی انتخابی در تشنج های ابسانس:'
$ws.Range("A110").Value = 'رفلکس‌های هرینگ-بروئر توسط چه چیزی واسطه می‌شوند؟'
$ws.Range("A111").Value = 'تب کیو چگونه منتقل می‌شود؟'
$ws.Range("A112").Value = 'ارزیابی کیفیت پروتئین به بهترین نحو توسط کدام روش انجام می‌شود؟'
$ws.Range("A113").Value = 'کدام یک از موارد زیر تست آگلوتیناسیون هتروفیل نیست؟'
$ws.Range("A114").Value = 'بیمار با تاول دردناک در طول دیواره قفسه سینه مراجعه می‌کند. کدام یک از آزمایش‌های زیر برای تشخیص مفید نیست؟'
$ws.Range("A115").Value = 'مهم‌ترین عامل خطر برای سندرم دیسترس تنفسی نوزادی چیست؟'
$ws.Range("A116").Value = 'شستشوی معده در کدام مورد منع مصرف دارد؟'
$ws.Range("A117").Value = 'همه موارد زیر روش‌های درمانی قطعی برای کارسینوم هپاتوسلولار هستند، به جز:'
$ws.Range("A118").Value = 'کدام یک از موارد زیر از ویژگی‌های نارسایی سمت راست قلب نیست؟'
$ws.Range("A119").Value = 'فرآیند توزیع مجدد با کدام یک مشاهده می‌شود؟'
$ws.Range("A120").Value = 'ظاهر چرخ مانند پرده صماخ در asom به دلیل چیست؟'
$ws.Range("A121").Value = 'یک مرد ۴۴ ساله با شروع ناگهانی درد شدید در ربع فوقانی راست شکم (ruq)، آسیت، هپاتومگالی حساس و هماتمز (استفراغ خونی) مراجعه می‌کند. این علائم نشان‌دهنده سندرم بود-کیاری است که ناشی از کدام مورد زیر می‌باشد؟'
$ws.Range("A122").Value = 'کدام گاز بیهوشی تمایل بالایی به تجمع در حفره‌ها دارد؟'
$ws.Range("A123").Value = 'در سیفلیس ثانویه همه موارد زیر دیده می‌شوند به جز -'
$ws.Range("A124").Value = 'کدام اصل قانون اساسی به ما اجازه می‌دهد افراد در اعتصاب غذا را مجبور به خوردن کنیم؟'
$ws.Range("A125").Value = 'ارتفاع رحم بیشتر از سن حاملگی بیمار در تمام موارد زیر به جز کدام یک است؟'
$ws.Range("A126").Value = 'تونسیلکتومی در کدام مورد نشان داده شده است؟'
$ws.Range("A127").Value = 'بر اساس قانون فرانک-استارلینگ، میزان پیش بار متناسب است با'
$ws.Range("A128").Value = 'کدام یک از موارد زیر در مورد عفونت سالمونلا صحیح است؟'
$ws.Range("A129").Value = 'حداقل سطح tsh که توسط آزمایش tsh نسل سوم تست‌های عملکرد تیروئید قابل تشخیص است، کدام یک از گزینه‌های زیر است؟'
$ws.Range("A130").Value = 'آغازگر برای کامپوزیت خوددرمانی چیست؟'
$ws.Range("A131").Value = 'کاتاراکت ثانویه/بعدی، کدورت عدسی است که پس از کدام مورد رخ می‌دهد؟'
$ws.Range("A132").Value = 'شایع ترین تومور استخوانی گوش میانی کدام است؟'
$ws.Range("A133").Value = 'درمان انتخابی برای یورتروسل'
$ws.Range("A134").Value = 'کدام یک از پارامترهای زیر برای ارزیابی خطر و نیاز به ضد انعقاد در بیمار مبتلا به فیبریلاسیون دهلیزی استفاده نمی‌شود؟'
$ws.Range("A135").Value = 'در انسان، جریان خون مؤثر کلیوی چقدر است؟'
$ws.Range("A136").Value = 'عدم جوش خوردن یک عارضه بسیار شایع در شکستگی‌های داخل کپسولی گردن فمور است. کدام یک از موارد زیر از علل بسیار مهم این عارضه محسوب نمی‌شود؟'
$ws.Range("A137").Value = 'بیماری هالبان به چه علتی رخ می‌دهد؟'
$ws.Range("A138").Value = 'آمفوتریسین b بر چه چیزی اثر می‌کند؟'
$ws.Range("A139").Value = 'همه موارد زیر از عناصر ارزیابی هستند به جز'
$ws.Range("A140").Value = 'کدام یک از موارد زیر از عوارض مرتبط با پاروتیدکتومی محسوب نمی‌شود؟'
$ws.Range("A141").Value = 'تمامی عبارات زیر در مورد دیورتیکول زنکر صحیح هستند به جز:'
$ws.Range("A142").Value = 'ظاهر چوب‌درام (drumstick) در کدام یک دیده می‌شود؟'
$ws.Range("A143").Value = 'کدام آنتی‌بادی برای تشخیص بیماری بافت همبند مختلط (mctd) استفاده می‌شود؟'
$ws.Range("A144").Value = 'حفره پر از مایع در مرکز نخاع در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A145").Value = 'این بیمار مبتلا به فشار خون بالا با فلج نیمه راست بدن بستری شده است. سی‌تی‌اسکن ساده چه چیزی را نشان می‌دهد؟'
$ws.Range("A146").Value = 'یک مرد 40 ساله، سیگاری مزمن، با ناراحتی حاد اپیگاستریک در طول یک ساعت گذشته مراجعه می‌کند. نوار قلب او افزایش قطعه st در لیدهای تحتانی را نشان می‌دهد. مداخله فوری چیست؟'
$ws.Range("A147").Value = 'سندرم بدخیم نورولپتیک ناشی از کدام دارو است؟ مارس 2009'
$ws.Range("A148").Value = 'روش بیلینگز برای پیشگیری از بارداری به ________ اشاره دارد.'
$ws.Range("A149").Value = 'یک پسر قصد دارد با اتوبوس سفر کند. کدام یک از داروهای زیر می‌تواند برای جلوگیری از بیماری حرکت در این فرد استفاده شود؟'
$ws.Range("A150").Value = 'در سال 1947 یک کمیته در ada به صورت نمادین کدام یک از سیستم های شماره گذاری دندان را توصیه کرد؟'
$ws.Range("A151").Value = 'یک مرد ۲۰ ساله در شهری دور از زادگاهش پیدا شد که در یک کارخانه کار می‌کرد. او قادر نبود درباره زندگی گذشته‌اش توضیح دهد. کدام یک از موارد زیر احتمالاً تشخیص صحیح است؟'
$ws.Range("A152").Value = 'در ph فیزیولوژیکی، مولکول‌های dna چگونه هستند؟'
$ws.Range("A153").Value = 'کودک با تورم خط میانی در گردن، محتمل‌ترین تشخیص چیست؟'
$ws.Range("A154").Value = 'یک مرد ۶۰ ساله می‌خواهد دندان مصنوعی برای او ساخته شود. او از داروهای خوراکی کاهش‌دهنده قند خون استفاده می‌کند. معاینه داخل دهانی، قوس‌های بی‌دندان با لبه‌های محکم و ارتفاع خوب آلوئولار را نشان می‌دهد. کدام یک از موارد زیر را می‌توان برای قالب‌گیری در این بیمار استفاده کرد؟'
$ws.Range("A155").Value = 'یک نیروی منفرد به تاج دندان‌های پیشین فک بالا اعمال می‌شود. چرخش حول نقطه چرخش موجود در کدام ناحیه اتفاق می‌افتد؟'
$ws.Range("A156").Value = 'یک پزشک باید یک عمل جراحی فوری روی یک بیمار بیهوش انجام دهد تا جان او را نجات دهد. اما هیچ یک از بستگان برای گرفتن رضایت حضور ندارند. او بدون گرفتن رضایت از任何人 عمل را انجام می‌دهد؛ او از اصل چه چیزی استفاده می‌کند؟ cmc (vellore) 09'
$ws.Range("A157").Value = 'عملکرد پلاکت را می توان با کدام روش ارزیابی کرد؟'
$ws.Range("A158").Value = 'قانون کورن مربوط به چیست؟'
$ws.Range("A159").Value = 'مزسان همراه با سیکلوفسفامید تجویز می‌شود تا؟'
$ws.Range("A160").Value = 'پاسخ برانگیختگی توسط چه چیزی واسطه می‌شود؟'
$ws.Range("A161").Value = 'نوریت اپتیک باعث چه چیزی می‌شود؟'
$ws.Range("A162").Value = 'بیمار تحت درمان با لیتیوم همچنین از فشار خون بالا رنج می‌برد. کدام یک از داروهای ضد فشار خون زیر در صورت تجویز همراه با لیتیوم می‌تواند باعث مسمومیت در این بیمار شود؟'
$ws.Range("A163").Value = 'شایع ترین تومور کبد کدام است؟'
$ws.Range("A164").Value = 'شایع‌ترین عارضه پس از برداشتن پروستات از طریق مجرای ادرار (turp) چیست؟'
$ws.Range("A165").Value = 'لیپودرماتواسکلروز بیشتر در کدام ناحیه دیده می‌شود؟'
$ws.Range("A166").Value = 'ضریب تنفسی کربوهیدرات چیست؟'
$ws.Range("A167").Value = 'شایع‌ترین علت کوری قشری چیست؟'
$ws.Range("A168").Value = 'کدام یک از پاسخ‌های بافتی میزبان زیر در التهاب حاد دیده نمی‌شود؟'
$ws.Range("A169").Value = 'متهموگلوبینمی با کدامیک از موارد زیر مرتبط است؟'
$ws.Range("A170").Value = 'چرخه میکروفیلاریا در پشه -'
$ws.Range("A171").Value = 'تمام موارد زیر در مورد حذف لایه اسمیر دنتین در طول آماده‌سازی کانال ریشه صحیح هستند به جز:'
$ws.Range("A172").Value = 'گنوکوک کدام ماده را تخمیر می‌کند؟'
$ws.Range("A173").Value = 'کدام یک از عبارات زیر در مورد عضله اسکلتی و قلبی نادرست است؟'
$ws.Range("A174").Value = 'در مرحله الیگوری نارسایی کلیه، همه موارد زیر صحیح هستند، به جز -'
$ws.Range("A175").Value = 'یک زن ۳۰ ساله g2p1 تحت عمل سزارین انتخابی تکرار در ترم قرار دارد. نوزاد بدون هیچ مشکلی به دنیا می‌آید، اما جفت به راحتی قابل خارج شدن نیست زیرا صفحه واضحی بین جفت و دیواره رحم شناسایی نمی‌شود. جفت به صورت تکه‌تکه خارج می‌شود. پس از آن، آتونی رحم و خونریزی رخ می‌دهد. تمام موارد زیر در مورد این وضعیت صحیح است به جز'
$ws.Range("A176").Value = 'فلورسکامین برای تشخیص چه چیزی استفاده می‌شود؟'
$ws.Range("A177").Value = 'مصرف الکل در دوران بارداری باعث کدام یک از موارد زیر نمی‌شود؟'
$ws.Range("A178").Value = 'یک بیمار 50 ساله سیگاری مزمن با سرفه، هموپتیزی (سرفه خونی) و کاهش وزن بدون دلیل مراجعه کرده است. سی‌تی اسکن یک تومور در ریه راست که به پلور مدیاستین حمله کرده است را نشان می‌دهد. بیمار احتمالاً درد ارجاعی را در کدام یک از نواحی زیر تجربه خواهد کرد؟'
$ws.Range("A179").Value = 'بالون والووتومی در تمام موارد زیر موفقیت‌آمیز است به جز –'
$ws.Range("A180").Value = 'آزمایش vdrl برای تشخیص سیفلیس نمونه ای از کدام تست است؟'
$ws.Range("A181").Value = 'کدام یک از موارد زیر به عنوان شتاب‌دهنده در سیمان اکسید روی اوژنول استفاده نمی‌شود؟'
$ws.Range("A182").Value = 'تکنیک ترجیحی برای زایمان بدون درد'
$ws.Range("A183").Value = 'کدام یک از موارد زیر بخشی از نمره child-pugh نیست؟'
$ws.Range("A184").Value = 'شایع‌ترین نوع اتیولوژیک یوئیت چیست؟'
$ws.Range("A185").Value = 'کدام یک از موارد زیر کتوز نیست؟'
$ws.Range("A186").Value = 'بورلیا به دلیل چه چیزی دچار تغییرات آنتی ژنی می‌شود؟'
$ws.Range("A187").Value = 'تمام موارد زیر در مورد esbl صحیح هستند به جز:'
$ws.Range("A188").Value = 'کلسیفیکاسیون غضروف منیسک از ویژگی های کدام بیماری است؟'
$ws.Range("A189").Value = 'یک مرد 60 ساله که یک سال پیش تحت عمل فیکوآمولسیفیکاسیون قرار گرفته است، با شکایت از تاری دید در آن چشم طی یک ماه گذشته مراجعه کرده است. شروع تاری دید تدریجی بوده و با درد، قرمزی یا اشکریزی همراه نبوده است. او می‌گوید پس از عمل، دیدش واضح بوده است. شایع‌ترین علت این وضعیت چیست؟'
$ws.Range("A190").Value = 'در مورد کارسینوم فیبرولاملار کبد کدام گزینه صحیح است؟'
$ws.Range("A191").Value = 'عضلات ذاتی زبان از کجا منشأ می‌گیرند؟'
$ws.Range("A192").Value = '"متاستاتیک" کلسیفیکیشن اغلب در کدام یک دیده می‌شود؟'
$ws.Range("A193").Value = 'موارد زیر همگی از ویژگی‌های هایپوترمی شدید هستند به جز-'
$ws.Range("A194").Value = 'روش ترشح غدد سباسه چیست؟'
$ws.Range("A195").Value = 'ظرفیت باقیمانده عملکردی چیست؟'
$ws.Range("A196").Value = 'همه موارد زیر از عوارض کولیت اولسراتیو هستند، به جز:'
$ws.Range("A197").Value = 'ناحیه بی‌دندانی دوطرفه در قسمت خلفی دندان‌های باقی‌مانده را می‌توان به عنوان کندی کلاس ... طبقه‌بندی کرد.'
$ws.Range("A198").Value = 'به طور معمول، دررفتگی تحتانی دو طرفه عدسی در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range("A199").Value = 'شایع‌ترین علت دفع چند قطره خون از مقعد در یک کودک نوپا چیست؟ (aiims may 2013, may 2012)'
$ws.Range("A200").Value = 'بخیه‌های جذب‌شدنی کدام‌اند؟'
$ws.Range("A201").Value = 'یک کودک 10 ساله با درد و توده در ناحیه کمری راست بدون تب، با خم شدن مفصل ران راست و تغییرات ستون فقرات در عکسبرداری اشعه ایکس. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A202").Value = 'پلاسمای منجمد تازه در چه مواردی استفاده می‌شود؟'
$ws.Range("A203").Value = 'کدام یک از داروهای زیر به احتمال زیاد در برابر سویه‌های مقاوم به چند داروی مایکوباکتریوم توبرکلوزیس، از جمله سویه‌های مقاوم به استرپتومایسین، مؤثر است؟'
$ws.Range("A204").Value = 'باقیمانده نوتوکورد کدام است؟'
$ws.Range("A205").Value = 'کدام یک از اعصاب زیر در شکستگی تنه استخوان بازو آسیب می‌بیند؟'
$ws.Range("A206").Value = 'یک مرد 70 ساله مبتلا به فشار خون سیستولیک ایزوله است. در معاینه، فشار خون او 170/80 mmhg است و قلب و ریه‌های او طبیعی هستند. او هیچ بیماری دیگری ندارد. برای این بیمار با فشار خون بالا، مناسب‌ترین دارو را انتخاب کنید.'
$ws.Range("A207").Value = 'شایع‌ترین بیوتایپ استافیلوکوکوس اورئوس که باعث عفونت انسانی می‌شود کدام است؟'
$ws.Range("A208").Value = 'یک نوزاد نارس با سن حاملگی 30 هفته و وزن 1.2 کیلوگرم، در عرض 2 ساعت پس از تولد دچار دیسترس تنفسی، خرخر و تورفتگی قفسه سینه شده است. تعداد تنفس او 70 در دقیقه است. مناسب‌ترین اقدام بعدی در مدیریت این نوزاد کدام است؟'
$ws.Range("A209").Value = 'در دوران بارداری حساسیت تنفسی به دی‌اکسید کربن به دلیل سطوح بالاتر گردش خون کدام هورمون افزایش می‌یابد؟'
$ws.Range("A210").Value = 'کدام یک از مراحل زیر مربوط به ظهور تفکر شهودی در نظریه ژان پیاژه است؟'
$ws.Range("A211").Value = 'هایپربیلی روبینمی غیرکونژوگه در تمام موارد زیر دیده می‌شود به جز'
$ws.Range("A212").Value = 'تومور ماتریکس ساز کدام است؟'
$ws.Range("A213").Value = 'کندیلار کارتیلاژ در فک پایین به عنوان چه نوع کارتیلاژی شناخته می‌شود؟'
$ws.Range("A214").Value = 'در بلوک عصبی محیطی، کدام یک از بیحس کننده های موضعی زیر طولانی ترین مدت اثر را دارد؟ (یا) کدام یک از عوامل بیحس کننده موضعی زیر برای جراحی های طولانی مدت ترجیح داده می شود؟'
$ws.Range("A215").Value = 'توصیف نمودار شجره‌نامه‌ای'
$ws.Range("A216").Value = 'داروی شیمی‌درمانی فلورواوراسیل در بدن دستخوش یک سری تغییرات شیمیایی می‌شود که منجر به تشکیل کمپلکس کووالانسی می‌گردد، به طوری که به هر دو آنزیم تیمیدیلات سنتاز و متیلن-تتراهیدروفولات متصل می‌شود. مهار تشکیل دئوکسی تیمیدیلات و در نتیجه مسدود شدن تقسیم سلولی ناشی از چیست؟'
$ws.Range("A217").Value = 'لیگامان اینگوینال به دلیل کدام مورد تشکیل می‌شود؟'
$ws.Range("A218").Value = 'رزت‌های یافت شده در حاشیه پلک، ویژگی کدام یک از موارد زیر هستند؟'
$ws.Range("A219").Value = 'پرتودرمانی به بافت‌های هیپوکسیک ممکن است با درمان با کدام یک از موارد زیر تقویت شود؟'
$ws.Range("A220").Value = 'سندرم رگرسیون کادال در نوزادان مادران مبتلا به کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A221").Value = 'کدام یک از شرایط زیر با افزایش خطر ابتلا به آدنوکارسینوم مری همراه است؟'
$ws.Range("A222").Value = 'کدام یک از موارد زیر یک عامل ایمنی‌تحریک‌کننده است؟'
$ws.Range("A223").Value = 'استرپتوکوک و پنوموکوک توسط کدام مورد از هم تفکیک می‌شوند؟  
الف) حلالیت در نمک صفراوی  
ب) ویژگی‌های رشد  
ج) رنگ‌آمیزی گرم  
د) تخمیر اینولین'
$ws.Range("A224").Value = 'کدام یک از موارد زیر یک تعدیل کننده انتخابی گیرنده استروژن است؟'
$ws.Range("A225").Value = 'عملکرد پراکسیزوم ها چیست؟'
$ws.Range("A226").Value = 'در 100 زن که به مدت 20 ماه با استفاده از قرص‌های ضدبارداری (ocps) تحت پیگیری قرار گرفتند، 5 مورد باردار شدند. شاخص پرل را محاسبه کنید.'
$ws.Range("A227").Value = 'تغییر مشاهده شده در نوار قلب در هایپرکالمی چیست؟'
$ws.Range("A228").Value = 'یک مرد ۶۰ ساله با علائم پلی اوری و شب ادراری مراجعه می‌کند. بررسی‌ها نشان‌دهنده سطح نرمال گلوکز ناشتا و کلسیم است. الکترولیت‌های ادراری او نشان‌دهنده نقص در انتقال کلیوی است. کدام یک از موارد زیر محتمل‌ترین علت علائم اوست؟'
$ws.Range("A229").Value = 'کدام یک از ویژگی‌های هموگلوبینوری حمله‌ای شبانه نیست؟'
$ws.Range("A230").Value = 'تمامی موارد زیر در مورد مجرای شریانی باز (pda) صحیح هستند به جز –  
الف) بسته شدن خودبه‌خودی در برخی نوزادان ترم اتفاق می‌افتد  
ب) فشار خون ریوی ایجاد می‌شود  
ج) اندوکاردیت باکتریایی در pdaهای کوچک شایع‌تر است  
د) عفونت‌های مکرر قفسه سینه و نارسایی احتقانی ممکن است ایجاد شود  
ه) وجود آناتومیک pda نشانه‌ای برای جراحی است'
$ws.Range("A231").Value = 'آمنوره اولیه در همه موارد زیر دیده می‌شود به جز:'
$ws.Range("A232").Value = 'در مورد rhinosporidium seebri کدام گزینه صحیح است؟'
$ws.Range("A233").Value = 'در مورد ماموگرافی غربالگری کدام گزینه صحیح است؟'
$ws.Range("A234").Value = 'اصطلاح "آپاندیسیت سمت چپ" که به طور رایج استفاده می‌شود، در واقع چیست؟'
$ws.Range("A235").Value = 'در مورد راشیتیسم ناشی از کمبود ویتامین d کدام گزینه صحیح است؟'
$ws.Range("A236").Value = 'علامت "بنت اینر تیوب" در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A237").Value = 'پوشش های پلاستیکی سرنگ ها در کدام کیسه دفع می شوند؟'
$ws.Range("A238").Value = 'یک زن 60 ساله پس از زمین خوردن، اندام تحتانی او در حالت کشیده و چرخش به بیرون قرار دارد. تشخیص محتمل چیست؟'
$ws.Range("A239").Value = 'یک بیمار مبتلا به مننژیت سلی تشخیص داده شده و درمان ضد سل برای او آغاز شده است. با توجه به تداوم ضعف صورت علیرغم شروع درمان ضد سل، پزشک تصمیم به تجویز دوز بالای استروئید می‌گیرد. در طی یک هفته پس از شروع پردنیزولون، بیمار درد تابنده به پشت پای راست و مشکل در دورسی فلکشن پای راست را تجربه می‌کند. این علائم به دلیل کدام مورد است؟'
$ws.Range("A240").Value = 'شایع‌ترین درگیری چشمی در سارکوئیدوز چیست؟'
$ws.Range("A241").Value = 'پرالیدوکسیم در مسمومیت با کدام یک مؤثر نیست:'
$ws.Range("A242").Value = 'قانون ختم بارداری پزشکی سال ۱۹۷۱ در سال ۲۰۰۲ اصلاح شد تا شامل موارد زیر شود'
$ws.Range("A243").Value = 'در پلاسما، اگر ph برابر با 5 باشد، نسبت باز به اسید چقدر است؟'
$ws.Range("A244").Value = 'چه کسی مفهوم تولیدمندی در مقابل رکود را ارائه داد؟'
$ws.Range("A245").Value = 'نارسایی مزمن آدرنال با کدام یک از موارد زیر مشخص می‌شود؟'
$ws.Range("A246").Value = 'در سنتز پروتئین یوکاریوتی، فسفریلاسیون کدام یک از موارد زیر توسط انسولین افزایش می‌یابد؟'
$ws.Range("A247").Value = 'روش بالینی خوب (gcp) در تمام موارد زیر دیده می‌شود به جز؟'
$ws.Range("A248").Value = 'موارد استفاده سیکلوسپورین همه موارد زیر است به جز -'
$ws.Range("A249").Value = 'غدد زایس عبارتند از:'
$ws.Range("A250").Value = 'کدام یک از موارد زیر از ویژگی‌های سندرم پیر-روبین نیست؟'
$ws.Range("A251").Value = 'کدام یک از موارد زیر در مورد کارسینوم فولیکولار تیروئید صحیح است؟'
$ws.Range("A252").Value = 'سلول خاردار در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A253").Value = 'نوع پیوندی که برای پیوند کلیه مناسب‌تر است؟'
$ws.Range("A254").Value = 'کدام یک از موارد زیر در استراتژی ''safe'' گنجانده نشده است؟'
$ws.Range("A255").Value = 'فردی فکر می‌کند که دیگران یک ماشین روی سر او کار گذاشته‌اند. این مثال چیست؟'
$ws.Range("A256").Value = 'تغییر رنگ پوست به آبی-خاکستری'
$ws.Range("A257").Value = 'مفصل شانه از نظر آناتومیکی در کدام ناحیه ضعیف است؟'
$ws.Range("A258").Value = 'طبق قانون esi، مزایای زایمان در صورت سقط جنین برای چند روز مجاز است:'
$ws.Range("A259").Value = 'یک دختر ۱۶ ساله با آمنوره اولیه مراجعه می‌کند. ویژگی‌های جنسی ثانویه او به خوبی توسعه یافته‌اند و موهای زیر بغل و ناحیه تناسلی خوبی دارد. معاینه دستگاه تناسلی خارجی نشان‌دهنده عدم وجود واژن است. سونوگرافی عدم وجود رحم را نشان می‌دهد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A260").Value = 'قرص های ضد بارداری خوراکی در همه موارد زیر منع مصرف دارند به جز:'
$ws.Range("A261").Value = 'یک کشور در صورتی گواهی ریشه‌کنی فلج اطفال را دریافت می‌کند که در طول سال‌های اخیر هیچ مورد تأیید شده ویروس‌شناسی از فلج اطفال تشخیص داده نشده باشد -'
$ws.Range("A262").Value = 'در مورد پایش سنتینل کدام گزینه صحیح است؟ مارس 2011'
$ws.Range("A263").Value = 'همه موارد زیر در مورد اختلال تبدیلی درست است به جز:'
$ws.Range("A264").Value = 'ظاهر "پیچ بازکن" در مری در کدام شرایط مشاهده می‌شود؟'
$ws.Range("A265").Value = 'کدام یک از موارد زیر کم‌ترین محل شایع برای ضایعات آترواسکلروتیک است؟'
$ws.Range("A266").Value = 'کدام یک از عبارات زیر در مورد سندرم فون هیپل-لینداو صحیح است؟'
$ws.Range("A267").Value = 'در کدام یک از شرایط زیر می‌توان از ید رادیواکتیو در بیماری گریوز استفاده کرد؟'
$ws.Range("A268").Value = 'خطوط بلاشکو به خطوط کدام یک از موارد زیر اشاره دارد؟'
$ws.Range("A269").Value = 'زوال عقل برگشت‌پذیر از ویژگی‌های کدام یک از موارد زیر است؟'
$ws.Range("A270").Value = 'در مورد موکورمایکوزیس کدام گزینه صحیح است؟'
$ws.Range("A271").Value = 'قدرت یک چشم کاهش یافته به طور معمول چقدر است؟'
$ws.Range("A272").Value = 'در مورد اپیدمیولوژی آنفلوانزا کدام گزینه صحیح است؟'
$ws.Range("A273").Value = 'ویژگی مشخص آپوپتوز در میکروسکوپ نوری چیست؟'
$ws.Range("A274").Value = 'همه موارد زیر در مورد ادم ایدیوپاتیک زنان صحیح است به جز -'
$ws.Range("A275").Value = 'صدای اول قلب در چه مرحله‌ای رخ می‌دهد؟'
$ws.Range("A276").Value = 'کدام دارو به شدت وستیبولوتوکسیک است؟'
$ws.Range("A277").Value = 'کدام یک از گزینه‌های زیر درمان انتخابی برای یک کودک 2 ماهه با علائم اشک‌ریزش و ترشح از چشم با یافته‌های طبیعی اتاق قدامی و خلفی است؟'
$ws.Range("A278").Value = 'یک فرد 33 ساله الکلی که تحت درمان att است، با افزایش آهن سرم و افزایش اشباع ترانسفرین مراجعه می‌کند. تشخیص-'
$ws.Range("A279").Value = 'یک مرد 55 ساله که در 10 روز گذشته استراحت مطلق داشته است، از تنگی نفس و درد قفسه سینه شکایت دارد. عکس قفسه سینه طبیعی است. بررسی بعدی باید چیست؟'
$ws.Range("A280").Value = 'در صورت مرگ ناگهانی یک فوتبالیست جوان، اولین شک بالینی بر کدام یک از موارد زیر قرار می‌گیرد؟'
$ws.Range("A281").Value = 'یک بیمار مرد 30 ساله دارای زخم بزرگ، گسترده و پررنگ با بافت گرانولاسیون قرمز روشن روی گلانس آلت تناسلی است. هیچ لنفادنوپاتی مشاهده نشد. محتمل‌ترین عامل ایجادکننده این بیماری چیست؟'
$ws.Range("A282").Value = 'رادیوفرکوئنسی ابلیشن برای کدام مورد انجام می‌شود؟'
$ws.Range("A283").Value = 'کدام یک از موارد زیر ویژگی مشخصه صرع میوکلونیک نوجوانی است؟'
$ws.Range("A284").Value = 'کشف جان اسنو در مورد وبا نمونه ای از کدام مورد است؟'
$ws.Range("A285").Value = 'مدت زمان لازم برای حرکت دندان پرمولر در حال رویش به اندازه 1 میلی متر در استخوان، که در رادیوگرافی بایت-وینگ اندازه گیری می شود، چقدر است؟'
$ws.Range("A286").Value = 'اولین خط دارویی برای پوکی استخوان پس از یائسگی چیست؟'
$ws.Range("A287").Value = 'کدام یک از موارد زیر بخشی از مدیریت گازگرفتگی درجه سه سگ آلوده به هاری نیست: مارس 2009، سپتامبر 2012'
$ws.Range("A288").Value = 'نوع کم‌شنوایی مشاهده‌شده در مراحل اولیه بیماری منیر:'
$ws.Range("A289").Value = 'رفلکس مورو در ________ از بین می‌رود.'
$ws.Range("A290").Value = 'روبئوزیس ایریدیس در کدام یک از موارد زیر دیده نمی‌شود؟ سپتامبر 2008'
$ws.Range("A291").Value = 'تشخیص بیماری هوچکین بر اساس تمام موارد زیر است، به جز:'
$ws.Range("A292").Value = 'ویروس hiv در چه سالی کشف شد؟'
$ws.Range("A293").Value = 'انتقال دو ماده در یک جهت مشترک چه نامیده می‌شود؟'
$ws.Range("A294").Value = 'شاخص درمانی دارو را از اطلاعات داده شده در نمودار زیر پیدا کنید'
$ws.Range("A295").Value = 'در مورد فلج بل کدام گزینه صحیح است؟'
$ws.Range("A296").Value = 'از دست دادن ارتباط بین فردی چگونه مدیریت می‌شود؟'
$ws.Range("A297").Value = 'کدام گزینه در مورد hcg نادرست است؟'
$ws.Range("A298").Value = 'آمین های وازواکتیو که در التهاب نقش دارند کدامند؟'
$ws.Range("A299").Value = 'نوزادی با سابقه عدم دفع مکونیوم آورده شده است و در معاینه هیچ سوراخ مقعدی مشاهده نمی‌شود، فقط یک فرورفتگی وجود دارد. بررسی انتخابی چیست؟'
$ws.Range("A300").Value = 'بهترین پیش‌آگهی برای کارسینوم پانکراس در کدام ناحیه است؟'
$ws.Range("A301").Value = 'خط شنتون در عکس رادیوگرافی کدام ناحیه دیده می‌شود؟'
$ws.Range("A302").Value = 'یک مرد 45 ساله با یک توده بدون درد در گردن مراجعه می‌کند. یک تومور سفت و متحرک به اندازه 4 سانتی‌متر در زاویه فک چپ شناسایی می‌شود. بیوپسی تومور نشان‌دهنده سلول‌های میواپیتلیال در هم آمیخته با مناطق میگزوئید، موکوئید و غضروفی است. تومور به صورت جراحی برداشته می‌شود. کدام یک از موارد زیر محتمل‌ترین پیش‌آگهی است؟'
$ws.Range("A303").Value = 'مسیر پارووسلولار، از هسته زانویی جانبی به قشر بینایی، سیگنال‌های مربوط به تشخیص چه چیزی را حمل می‌کند؟'
$ws.Range("A304").Value = 'تمامی عبارات زیر در مورد وازدفران صحیح هستند، به جز:'
$ws.Range("A305").Value = 'در مورد بیضه نزول نکرده، پس از چه سنی دیگر نزول اتفاق نمی‌افتد؟'
$ws.Range("A306").Value = 'کدام دارو در نارسایی کلیوی استفاده نمی‌شود؟'
$ws.Range("A307").Value = 'یک مرد 18 ساله به نام روهیت از ضعف شدید، سرگیجه و خواب‌آلودگی در طول 4 هفته گذشته شکایت داشت. در یک دوره قبلی 3 هفته پیش، او دچار خواب‌آلودگی و خستگی عمومی شده بود. در هنگام پرسش، او تشنگی غیرقابل کنترل و نیاز مکرر به ادرار کردن را گزارش داد. او 4 کیلوگرم وزن از دست داده بود. آزمایش‌های انجام شده در زمان پذیرش، سطح گلوکز 560 میلی‌گرم در دسی‌لیتر را نشان داد. نمونه ادرار 4+ برای گلوکز و "مقدار زیاد" استون داشت. hba1c او 14% بود. کدام یک از موارد زیر در مورد این بیمار صحیح است؟'
$ws.Range("A308").Value = 'نوتوکورد به چه چیزی تبدیل می‌شود؟'
$ws.Range("A309").Value = 'عوامل اپی ژنتیک کنترل کننده رشد اسکلت کدامند؟'
$ws.Range("A310").Value = 'در hiv، کدام بخش از مغز دچار آتروفی می‌شود؟'
$ws.Range("A311").Value = 'لکه‌های کوپلیک نشانه‌ی پاتوگنومونیک کدام بیماری هستند؟'
$ws.Range("A312").Value = 'نوع آمبولی که در شکستگی استخوان‌های بلند دیده می‌شود:'
$ws.Range("A313").Value = 'همی پلژی اغلب ناشی از ترومبوز کدام یک از موارد زیر است؟'
$ws.Range("A314").Value = 'کاهش قلیاییت ca(oh)2 هنگامی که به عنوان پانسمان کانال ریشه استفاده می‌شود، از چه زمانی شروع می‌شود؟'
$ws.Range("A315").Value = 'در مورد کلامیدیا همه موارد زیر صحیح هستند به جز:'
$ws.Range("A316").Value = 'استرپتوکوکال tss ناشی از چیست؟'
$ws.Range("A317").Value = 'طول اسپرماتوزوئید بالغ انسان چقدر است؟'
$ws.Range("A318").Value = 'رایج‌ترین آزمایش ترپونمایی مورد استفاده در تشخیص سیفلیس کدام است؟'
$ws.Range("A319").Value = 'همه موارد زیر از علائم مرگ هستند، به جز:'
$ws.Range("A320").Value = 'در پسران اولین نشانه قابل مشاهده بلوغ چیست؟'
$ws.Range("A321").Value = 'یک زن ۳۵ ساله با علائم پارستزی و ضعف در اندام تحتانی دو طرفه، احساس نواری شکل سفتی در اطراف تنه، کاهش دید دردناک در هر دو چشم همراه با دیپلوپی و درد پری اوربیتال به درمانگاه پزشکی مراجعه کرده است. سابقه حملات مشابه در گذشته با دوره‌های بهبودی بین حملات وجود دارد. در معاینه: آتاکسی وجود دارد، پاپیلیت (در بررسی فوندوس)، میوکیمیای صورت، بی‌اختیاری ادرار و یبوست مشاهده می‌شود. بررسی مایع مغزی-نخاعی پلئوسیتوز مونونوکلئر و افزایش igg را نشان می‌دهد. کدام یک از داروهای خوراکی زیر برای این شرایط تأیید شده‌اند؟  
۱. فینگولیمود  
۲. ناتالیزوماب  
۳. تری فلونوماید  
۴. گلاتیرامر استات'
$ws.Range("A322").Value = 'عدم تحمل ارثی فروکتوز به دلیل کمبود کدام آنزیم است؟'
$ws.Range("A323").Value = 'تظاهرات عصبی فلوئوروزیس زمانی رخ می‌دهد که مدت زمان مصرف آن چند سال باشد؟'
$ws.Range("A324").Value = 'کدام یک از مقادیر زیر نشان‌دهنده نکروز حاد توبولی است؟'
$ws.Range("A325").Value = 'چه درمانی نباید برای بیمار با درگیری فورکیشن درجه ii مورد نیاز باشد؟'
$ws.Range("A326").Value = '«پروتئین‌های جی» چه هستند؟'
$ws.Range("A327").Value = 'یک داروی آزمایشی ترجیحاً شریانچه‌های خودتنظیم‌کننده کرونری را گشاد می‌کند بدون آنکه بر رگ‌های هدایت‌کننده بزرگ تأثیر بگذارد. این دارو به احتمال زیاد:'
$ws.Range("A328").Value = 'کدام ناحیه تالاموس در انسفالوپاتی ورنیکه درگیر می‌شود؟'
$ws.Range("A329").Value = 'بهترین روش برای ترمیم جراحی آسیب به آئورت شکمی بالاتر از سطح شریان کلیوی که شامل شریان فوقانی، تنه سلیاک و شریان فوق کلیوی است:'
$ws.Range("A330").Value = 'بیمار به مدت 3 روز از مصرف الکل خودداری کرده است. بیمار به مدت 20 سال گذشته الکلی شناخته شده است. ناگهان دچار لرزش بی‌قرار، آشفتگی می‌شود. تشخیص -'
$ws.Range("A331").Value = 'لامبدا محل تلاقی کدام است؟'
$ws.Range("A332").Value = 'یک مرد 65 ساله با سابقه انفارکتوس میوکارد برای کلهسیستکتومی انتخابی برنامهریزی شده است. وزن او 73 کیلوگرم و قد او 6 فوت و 2 اینچ است. حجم سلولهای فشرده خون او 35 است. حداکثر خونریزی مجاز برای این عمل چقدر است؟'
$ws.Range("A333").Value = 'رنگ آمیزی نقره متنامین گوموری برای قارچ در زیر نشان داده شده است. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A334").Value = 'یکی از ویژگی‌های واسکولیت کلیوی در کودکان چیست؟'
$ws.Range("A335").Value = 'پو د اورانژ به دلیل چیست؟'
$ws.Range("A336").Value = 'اندیکاسیون‌های جراحی در واریکوسل عبارتند از:'
$ws.Range("A337").Value = 'همه موارد زیر از موارد کاربرد مانیتورینگ پتانسیل برانگیخته هستند به جز:'
$ws.Range("A338").Value = 'کدام انگشت زودتر درگیر انقباض دوپویترن می‌شود؟'
$ws.Range("A339").Value = 'کدام یک از موارد زیر از مشتقات بافت همبند نیست؟'
$ws.Range("A340").Value = 'احتمال ابتلا به فیبروز کیستیک در صورتی که تنها یکی از والدین مبتلا باشد و دیگری سالم باشد'
$ws.Range("A341").Value = 'شایع ترین ناهنجاری مرتبط با سندرم دیسترس تنفسی حاد (ards) کدام است؟'
$ws.Range("A342").Value = 'کدام یک از ویژگی‌های زیر امکان استفاده از سلکوکسیب را در بیماران مبتلا به آرتریت فراهم می‌کند؟'
$ws.Range("A343").Value = 'کدام یک از اختلالات زیر ممکن است حتی در صورت شروع پس از 18 سالگی تشخیص داده شود؟'
$ws.Range("A344").Value = 'یافتن یک شریان نافی منفرد در سونوگرافی نشانگر چیست؟'
$ws.Range("A345").Value = 'یک زن پس از عمل کوله‌سیستکتومی باز تحت بیهوشی عمومی، پس از 12 ساعت دچار تنگی نفس می‌شود. کریپیتاسیون در لوب تحتانی راست وجود دارد. علت احتمالی چیست؟'
$ws.Range("A346").Value = 'کدام یک از موارد زیر در فشار خون حاد هنگام زایمان قابل تجویز نیست؟'
$ws.Range("A347").Value = 'شایع‌ترین علت خونریزی بین بطنی چیست؟'
$ws.Range("A348").Value = 'کدام یک از شل‌کننده‌های عضلانی اسکلتی زیر باعث درد در هنگام تزریق می‌شود؟'
$ws.Range("A349").Value = 'همه موارد به جز یکی برای استفاده در پیشگیری اورژانسی از بارداری توصیه می‌شوند'
$ws.Range("A350").Value = 'اوروگرانولوسیت ها عبارتند از:'
$ws.Range("A351").Value = 'بیماران با دررفتگی مادرزادی دوطرفه لگن (cdh) با کدام الگوی راه رفتن زیر راه می‌روند؟'
$ws.Range("A352").Value = 'فاصله زمانی داده شده در قطار ۴ چیست؟'
$ws.Range("A353").Value = 'متهم به تجاوز جنسی را می‌توان حتی بدون رضایت او تحت کدام بخش از آیین‌نامه‌ کیفری (crpc) معاینه پزشکی کرد؟'
$ws.Range("A354").Value = 'ابتکار ملی آهن پلاس شامل موارد زیر است:'
$ws.Range("A355").Value = 'عمل لوفور در چه مواردی انجام می‌شود؟'
$ws.Range("A356").Value = 'تبدیل محیطی پیشسازهای استروژن در بیماران چاق قبل از یائسگی منجر به تشکیل کدام ماده میشود؟'
$ws.Range("A357").Value = 'درمان انتخابی برای نوروبلاستوم مرحله 4s چیست؟'
$ws.Range("A358").Value = 'عبارت صحیح در مورد ترمیم برشی باز چیست؟'
$ws.Range("A359").Value = 'کدام یک از موارد زیر یک عفونت باکتریایی است؟'
$ws.Range("A360").Value = 'همودیالیز در کدام مورد مفید نیست؟'
$ws.Range("A361").Value = 'کدام مورد جزء سه‌گانه کارسینوم سلول کلیوی (renal cell carcinoma) نیست؟'
$ws.Range("A362").Value = 'رادیکال‌های آزاد می‌توانند توسط آنزیم‌های زیر غیرفعال شوند به جز:'
$ws.Range("A363").Value = 'رسوبات مزانژیال iga در گلومرول های کلیوی ویژگی کدام یک از موارد زیر است؟'
$ws.Range("A364").Value = 'کدام یک از ویژگی‌های زیر می‌تواند آدنوم پاراتیروئید را از هایپرپلازی متمایز کند؟'
$ws.Range("A365").Value = 'مرگ ناگهانی قلبی ممکن است در تمام موارد زیر رخ دهد، به جز:'
$ws.Range("A366").Value = 'یک زن 28 ساله به نام کاترینا با شکایت از درد سوپراپوبیک، تکرر ادرار و دیزوری به درمانگاه شما مراجعه می‌کند. او همچنین حدود 30 دقیقه پیش در آخرین ادرار خود خون دیده است. آزمایش ادرار او پیوری نشان می‌دهد اما هیچ کست گلبول سفیدی وجود ندارد. در معاینه فیزیکی، حساسیت سوپراپوبیک در لمس مشاهده می‌شود. کدام یک از موارد زیر احتمالاً تشخیص این بیمار است؟'
$ws.Range("A367").Value = 'کدام یک از عوامل خطر قابل اصلاح برای تصلب شرایین نیست؟'
$ws.Range("A368").Value = 'ساختار چهارم پروتئین چیست؟'
$ws.Range("A369").Value = 'بهترین شاخص برای محاسبه ارزش تغذیه‌ای پروتئین کدام یک از موارد زیر است؟'
$ws.Range("A370").Value = 'عکس رادیوگرافی زانو را به دقت بررسی کنید. تشخیص:'
$ws.Range("A371").Value = 'کدام مورد از ویژگی‌های سندرم هورنر نیست؟'
$ws.Range("A372").Value = 'حلقه والدایر شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A373").Value = 'یک بیمار 6 ساعت پس از گزیده شدن توسط مار با ادم موضعی خفیف در محل آسیب به اورژانس مراجعه می‌کند. در معاینه هیچ ناهنجاری مشاهده نمی‌شود و نتایج آزمایشگاهی طبیعی هستند. مناسب‌ترین اقدام کدام است؟'
$ws.Range("A374").Value = 'مقدار بحرانی شاخص شوک که در بارداری غیرطبیعی در نظر گرفته می‌شود چیست؟'
$ws.Range("A375").Value = 'گونوره (سوزاک) — کدام یک از موارد زیر از ویژگی‌های بالینی آن نیست؟'
$ws.Range("A376").Value = 'در دررفتگی قدامی شانه، کدام عصب معمولاً تحت تأثیر قرار می‌گیرد؟'
$ws.Range("A377").Value = 'اپی فارینکس همچنین نامیده می‌شود'
$ws.Range("A378").Value = 'تینه آ "اینگونیتو" به دلیل استفاده نامناسب از کدام یک از موارد زیر به وجود می‌آید؟'
$ws.Range("A379").Value = 'یک مرد 69 ساله یک مولتی ویتامین حاوی ویتامین b12 دریافت می‌کند. او احساس خوبی دارد، هیچ علائمی گزارش نمی‌دهد و معاینه فیزیکی او طبیعی است. آزمایش cbc او کاملاً طبیعی است. کدام یک از موارد زیر بیشترین ویژگی جذب ویتامین b12 را دارد؟'
$ws.Range("A380").Value = 'همه موارد زیر از سوراخ آئورتیک دیافراگم عبور می‌کنند، به جز:'
$ws.Range("A381").Value = 'در زنان باردار دارای آنتی بادی های تیروئید پراکسیداز، خطر افزایش یافته چه چیزی وجود دارد؟'
$ws.Range("A382").Value = 'تمامی موارد زیر در مورد بوتولیسم صحیح هستند به جز -'
$ws.Range("A383").Value = '''دادگاه نوجوانان'' تحت نظارت چه کسی است؟ dnb 08'
$ws.Range("A384").Value = 'شاخص سوءتغذیه متوسط تا شدید کدام است؟'
$ws.Range("A385").Value = 'دیسمنوره احتقانی در کدام یک از بیماران مشاهده می‌شود؟'
$ws.Range("A386").Value = 'دیستروفی عضلانی دوشن صحیح است؟'
$ws.Range("A387").Value = 'کودکی برای مدت طولانی با رژیم غذایی اصلی ذرت تغذیه شده است. کدام یک از ویتامین‌های زیر ممکن است در بدن او کمبود داشته باشد؟'
$ws.Range("A388").Value = 'اصطلاح تومور جهانی به کدام مورد اشاره دارد؟'
$ws.Range("A389").Value = 'نازک‌ترین استخوان بخشی از کدام استخوان است؟'
$ws.Range("A390").Value = 'بریگما نامی است که به محل اتصال کدام موارد داده می‌شود؟'
$ws.Range("A391").Value = 'کدام گزینه در مورد اشرشیاکلی انتروتوکسیژنیک صحیح است؟'
$ws.Range("A392").Value = 'شایع ترین ناهنجاری مادرزادی پانکراس کدام است؟'
$ws.Range("A393").Value = 'لنف دندان‌ها به تمامی گره‌های زیر تخلیه می‌شود به جز'
$ws.Range("A394").Value = 'هنگامی که رکتوم با هوا از طریق یک لوله مقعدی پر می‌شود، در صورت آپاندیسیت درد و حساسیت در فوسای ایلیاک راست رخ می‌دهد. این حالت به عنوان چه نامیده می‌شود؟'
$ws.Range("A395").Value = 'سیروز میکرونودولار معمولاً در همه موارد زیر دیده می‌شود به جز -'
$ws.Range("A396").Value = 'علامت آسپیتز مربوط به کدام مورد است؟'
$ws.Range("A397").Value = 'سطح بحرانی الکل در خون که بالاتر از آن اختلال در رانندگی ایجاد می‌شود کدام است؟'
$ws.Range("A398").Value = 'نقش اصلی 2،3 بیس فسفوگلیسرات در گلبول‌های قرمز چیست؟'
$ws.Range("A399").Value = 'معیارهای نوزاد "در معرض خطر" شامل موارد زیر است به جز –'
$ws.Range("A400").Value = 'شایع‌ترین الگوی آندومتر در خونریزی رحمی عملکردی:'
$ws.Range("A401").Value = 'سلول‌های پاپکورن در کدام مورد یافت می‌شوند؟'
$ws.Range("A402").Value = 'اکلوژن متعادل دو طرفه بر اساس مفهوم ارائه شده توسط چه کسی است؟'
$ws.Range("A403").Value = 'مردی با درد قفسه سینه مراجعه می‌کند. نوار قلب کاهش قطعه st در لیدهای v1-v4 را نشان می‌دهد. کدام یک از موارد زیر نباید تجویز شود؟'
$ws.Range("A404").Value = 'کدام یک از شاخص های پلاک زیر از عامل افشا کننده استفاده نمی کند؟'
$ws.Range("A405").Value = 'برای کدام یک از فلج‌های عصبی زیر، آزمون زیر تشخیص داده می‌شود؟'
$ws.Range("A406").Value = 'فلوراید، که در جمع‌آوری نمونه‌های خون استفاده می‌شود، آنزیم زیر را مهار می‌کند: wb 07; pgi 08; fmge 10, 11'
$ws.Range("A407").Value = 'عبارات زیر را در مورد خواجه‌ها بخوانید: ۱) آلت تناسلی و کیسه بیضه در هجره‌ها قبل از بلوغ برداشته می‌شود ۲) اندام تناسلی زنانه‌ها دست‌نخورده باقی می‌ماند'
$ws.Range("A408").Value = 'فرم کامل aria چیست؟'
$ws.Range("A409").Value = 'ظاهر طحال با رسوبات آمیلوئید در پالپ سفید طحال به عنوان چه چیزی شناخته می‌شود؟'
$ws.Range("A410").Value = 'بتا بلاکرها:'
$ws.Range("A411").Value = 'کدام یک از داروهای زیر مدت زمان درمان ارتودنسی را افزایش می‌دهد؟'
$ws.Range("A412").Value = 'شایع‌ترین توده بینی:'
$ws.Range("A413").Value = 'ماده نگهدارنده مورد استفاده برای نمونه های سم شناسی:'
$ws.Range("A414").Value = 'کدام یک از موارد زیر برای طبقه‌بندی آدنوکارسینومای ریوی استفاده می‌شود؟'
$ws.Range("A415").Value = 'bal در کدام مورد منع مصرف دارد؟'
$ws.Range("A416").Value = 'درمان انتخابی برای کارسینوم گلوت t1n0m0 چیست؟'
$ws.Range("A417").Value = 'اجسام دوله در کدام یک از موارد زیر دیده می‌شوند؟'
$ws.Range("A418").Value = 'شایع‌ترین علامت در حاملگی خارج رحمی -'
$ws.Range("A419").Value = 'کدام یک از موارد زیر یک قند غیرکاهنده است؟'
$ws.Range("A420").Value = 'کدام یک از ویژگی های زیر در روان‌پریشی (سایکوز) گنجانده نشده است؟'
$ws.Range("A421").Value = 'سارکوئیدوز'
$ws.Range("A422").Value = 'کدام مورد در مورد راشیتیسم مقاوم به ویتامین d صحیح است؟'
$ws.Range("A423").Value = 'یک نوزاد دختر با اندام تناسلی مبهم، دچار کمبود 21-هیدروکسیلاز از نوع نمک‌سوز است. کدام یک از کاریوتایپ‌های زیر را انتظار دارید پیدا کنید؟'
$ws.Range("A424").Value = 'عبارت نادرست در مورد دیسمنوره اسپاسمودیک کدام است؟'
$ws.Range("A425").Value = 'یک پسر ۴ ساله آفریقایی دچار تودهای با رشد سریع در سمت راست صورتش میشود. نمونهبرداری از این ضایعه بهصورت زیر است. این نئوپلاسم با جابجاییهای کروموزومی مرتبط است که شامل کدام مورد میشود؟'
$ws.Range("A426").Value = 'آزمایش تکان دادن سر برای تشخیص چه موردی انجام می‌شود؟'
$ws.Range("A427").Value = 'تومور مرتبط با گسترش خارج جمجمه‌ای'
$ws.Range("A428").Value = 'نوتروپنی در تمام موارد زیر دیده می‌شود به جز:'
$ws.Range("A429").Value = 'یک زن ۶۴ ساله به دنبال شروع حاد ضعف در سمت چپ بدن در بیمارستان بستری می‌شود. او دچار سکته مغزی بزرگ در نیمکره راست مغز شده و در بستر استراحت می‌کند. در روز پنجم بستری، اشباع اکسیژن او در هوای اتاق به ۹۰٪ کاهش می‌یابد. بیمار حال عمومی خوبی دارد، ضعف عصبی تغییری نکرده است، فشار خون ۱۳۰/۹۰ میلی‌متر جیوه و نبض ۸۰ در دقیقه است. معاینه قفسه سینه کاهش فرمنیتس، کدری در پرکاشن و عدم وجود صداهای تنفسی در قسمت پایینی ریه چپ را نشان می‌دهد. همچنین انحراف نای به سمت چپ مشاهده می‌شود. با توجه به یافته‌های غیرطبیعی معاینه ریوی در این بیمار، محتمل‌ترین تشخیص را انتخاب کنید.'
$ws.Range("A430").Value = 'مسیر نزولی حرکتی مسئول حرکات سریع و ماهرانه کدام است؟'
$ws.Range("A431").Value = 'کدام یک از موارد زیر توسط یک کودک 3 ساله قابل انجام نیست؟'
$ws.Range("A432").Value = 'یک مرد ۲۰ ساله با حملات مکرر استفراغ خونی مراجعه می‌کند. سابقه‌ای از زردی یا نارسایی کبدی وجود ندارد. در معاینه، یافته‌های مهم شامل بزرگی متوسط طحال و وجود واریس مری است. آسیت یا زخم معده وجود ندارد. آزمایش‌های عملکرد کبد طبیعی هستند. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A433").Value = 'محدوده کامل محصولات cox می‌تواند توسط کدام‌یک سنتز شود؟'
$ws.Range("A434").Value = 'آزمون هس (تورنیکه) ویژگی کدام یک از موارد زیر است؟'
$ws.Range("A435").Value = 'مدار تنفس خودبه‌خودی مورد استفاده در کودکان کدام است؟'
$ws.Range("A436").Value = 'ائوزینوفیلی در بافت نکروز شده به دلیل چیست؟'
$ws.Range("A437").Value = 'ویتامین d به مقدار کافی در کدام ماده یافت نمی‌شود؟ مارس 2012'
$ws.Range("A438").Value = 'در هیستری، مکانیسم دفاعی مشاهده شده که مربوط به روانپزشکی نیست کدام است؟'
$ws.Range("A439").Value = 'کدام واکسن به شرایط نگهداری بسیار دقیق‌تری نیاز دارد؟'
$ws.Range("A440").Value = 'کدام مورد از ویژگی‌های مرور سیستماتیک نیست؟'
$ws.Range("A441").Value = 'در یک بیمار با زخم سوختگی که به اپیدرم سطحی گسترش یافته و درم را درگیر نکرده است، همه موارد زیر به جز کدام یک مشاهده می‌شود؟'
$ws.Range("A442").Value = 'نیمه عمر آلبومین پلاسما تقریباً چقدر است؟'
$ws.Range("A443").Value = 'نوع لیزر مورد استفاده برای ایریدوتومی چیست؟'
$ws.Range("A444").Value = 'محیط لوفلر چیست؟'
$ws.Range("A445").Value = 'محیط لوفر چیست؟'
$ws.Range("A446").Value = 'فردی دچار بریدگی در لبه قدامی عضله ذوزنقه ای در ناحیه گردن می‌شود و پس از آن نقطه شانه (کتف) او افتادگی پیدا می‌کند و در ابداکشن کامل بازو کمی مشکل دارد. به نظر می‌رسد کدام عصب قطع شده است؟'
$ws.Range("A447").Value = 'آوران رفلکس کرماستریک شامل کدام مورد است؟'
$ws.Range("A448").Value = 'یک پسر 5 ساله با تشخیص کلستئاتوم جیب جمع شونده خلفی فوقانی مواجه شده است. تمام موارد زیر بخشی از مدیریت این بیماری محسوب می‌شوند، به جز:'
$ws.Range("A449").Value = 'برون ده قلبی در کدام حالت کاهش می یابد؟'
$ws.Range("A450").Value = 'کدام یک از موارد زیر در مورد کلشی سین نادرست است؟'
$ws.Range("A451").Value = 'سوبسترا برای گلوکونئوژنز کدام است؟'
$ws.Range("A452").Value = 'یک پسر ده ساله با تشنج به بخش اورژانس کودکان مراجعه می‌کند. فشار خون در اندام فوقانی 200/140 میلی‌متر جیوه اندازه‌گیری شده است. نبض فمورال قابل لمس نبود. محتمل‌ترین تشخیص در میان موارد زیر کدام است؟'
$ws.Range("A453").Value = 'کدام یک از عضلات مشخص شده در نمودار داده شده بیشترین کمک را به جلو آوردن فک می‌کند؟'
$ws.Range("A454").Value = 'کدام رگ بیشتر در خونریزی زیر سخت‌شامه‌ای درگیر می‌شود؟'
$ws.Range("A455").Value = 'بهترین روش برای تخمین میزان پروتئینوری در یک کودک مبتلا به سندرم نفروتیک کدام است؟'
$ws.Range("A456").Value = 'شایع‌ترین شکستگی در دوران کودکی کدام است؟'
$ws.Range("A457").Value = 'بیماری کینباک چیست؟'
$ws.Range("A458").Value = 'نمونه‌برداری سیتولوژیک واژینال برای بررسی تغییرات هورمونی بهتر است از کدام ناحیه انجام شود؟'
$ws.Range("A459").Value = 'یک بیمار 40 ساله تحت عمل کوله سیستکتومی باز قرار گرفته است. جراح عمل گزارش کرده که عمل بدون عارضه بوده است. در روز اول پس از عمل، 100 میلی لیتر صفرا از درن قرار داده شده در بستر کیسه صفرا خارج شده است. در معاینه، بیمار بدون تب و زردی است. شکم نرم است و صداهای روده طبیعی شنیده می شود. به عنوان پزشک معالج، بهترین توصیه شما چیست؟'
$ws.Range("A460").Value = 'محلول نشانگر در دستگاه هوروکس حاوی چیست؟'
$ws.Range("A461").Value = 'پس از تخمک‌گذاری، تخمک‌های آزاد شده و لایه‌های اطراف سلول‌های گرانولوزا به عنوان چه شناخته می‌شوند؟'
$ws.Range("A462").Value = 'در صورت بیمار 42 ساله ضایعه ای مشاهده شد که در نمودار زیر نشان داده شده است. کدام یک از گزینه های زیر مدیریت ایده آل برای این شرایط می باشد؟ *نمودار'
$ws.Range("A463").Value = 'مهارکننده‌های پمپ پروتون برای بیماری زخم معده باید چگونه مصرف شوند؟'
$ws.Range("A464").Value = 'پرومترین در درمان کدام مورد استفاده می‌شود؟'
$ws.Range("A465").Value = 'منحنی تفکیک اکسیژن-هموگلوبین به کدام سمت جابجا می‌شود؟'
$ws.Range("A466").Value = 'صفرای سنگ ساز دارای ویژگی های زیر است -'
$ws.Range("A467").Value = 'همه موارد زیر گلومرولونفریت غیر تکثیری هستند به جز -'
$ws.Range("A468").Value = 'جی‌آی‌سی با ویسکوزیته بالا در چه سالی کشف شد؟'
$ws.Range("A469").Value = 'کلیه شامل چه ساختارهایی است؟'
$ws.Range("A470").Value = 'کدام rna دارای بالاترین درصد بازهای اصلاح شده است؟'
$ws.Range("A471").Value = 'لرزش سیستولیک در فضای بین دنده‌ای دوم یا سوم چپ در کدام یک از شرایط زیر شنیده می‌شود؟'
$ws.Range("A472").Value = 'تشنج‌های ناشی از سیپروفلوکساسین زمانی تشدید می‌شوند که کدام دارو به صورت همزمان تجویز شود؟'
$ws.Range("A473").Value = 'وظایف اصلی سلول های "i" چیست؟'
$ws.Range("A474").Value = 'اسپلنکتومی (برداشتن طحال) به عنوان روش درمان در کدام یک از شرایط زیر مفید نیست؟'
$ws.Range("A475").Value = 'عامل کاهش دهنده فشار داخل چشم که باعث آپنه در نوزادان می‌شود کدام است؟'
$ws.Range("A476").Value = 'در نمره بلوغ اسکلتی mckern و stewa، کدام یک از اتصالات اپی فیز زیر استفاده نمی‌شود؟'
$ws.Range("A477").Value = 'کدام یک از موارد زیر منبع تغذیه سلول های سرطانی است؟'
$ws.Range("A478").Value = 'رنگ سبز در سیستم تریاژ نشان دهنده چه چیزی است؟'
$ws.Range("A479").Value = 'استئوژنز ایمپرفکتا با دندان‌های طبیعی و صلبیه طبیعی:'
$ws.Range("A480").Value = 'یک مرد با علائم خستگی مراجعه می‌کند. آزمایش هماتولوژی نشان‌دهنده hb پایین و mcv بالا است. بررسی بعدی چیست؟'
$ws.Range("A481").Value = 'یک زن باردار مبتلا به سوزاک با ترشحات سفید. درمان انتخابی چیست؟'
$ws.Range("A482").Value = 'ارتباط دهنده بین کیسه بزرگ و کیسه کوچک-'
$ws.Range("A483").Value = 'بهترین درمان برای مسمومیت شدید با دیژیتال چیست؟'
$ws.Range("A484").Value = 'حرکت جانبی پیشرونده مفصل تمپورومندیبولار که با نرخی متناسب با حرکت رو به جلو اتفاق می‌افتد در کدام بخش رخ می‌دهد؟'
$ws.Range("A485").Value = 'کدام گزینه در مورد انتقال مواد در عرض غشای سلول صحیح است؟'
$ws.Range("A486").Value = 'همه موارد زیر در مولتیپل میلوما دیده می‌شوند، به جز -'
$ws.Range("A487").Value = 'آنتی بادی های ضد فودرین در کدام مورد مثبت هستند؟'
$ws.Range("A488").Value = 'استفاده از وادینگ در تفنگ لوله صاف همه موارد زیر را ایجاد می‌کند به جز -'
$ws.Range("A489").Value = 'مواد غذایی را به ترتیب نزولی میزان انرژی تولیدی آنها در هر 100 گرم رتبه بندی کنید. کربوهیدرات – a، چربی – b، الکل – c'
$ws.Range("A490").Value = 'مانیتورهای مدرن برای اندازه گیری etc02 از چه فناوری استفاده می کنند؟'
$ws.Range("A491").Value = 'nad+ توسط تمام آنزیم های زیر کاهش می یابد به جز:'
$ws.Range("A492").Value = 'لنفوم غیر هوچکینی با درجه پایین کدام است؟'
$ws.Range("A493").Value = 'اسپرمین در مایع منی توسط کدام آزمایش تشخیص داده می‌شود؟'
$ws.Range("A494").Value = 'اسپلنکتومی (برداشتن طحال) اغلب در مدیریت کدام یک از موارد زیر نشان داده می‌شود؟'
$ws.Range("A495").Value = 'بخش 376 ab قانون کیفری هند به چه موردی می‌پردازد؟'
$ws.Range("A496").Value = 'تغییر ژنتیکی در کدام یک از کروموزوم‌های زیر می‌تواند منجر به الیگودندروگلیوما شود؟'
$ws.Range("A497").Value = 'یک زن ۴۵ ساله با تنگی نفس شدید در بیمارستان بستری می‌شود. بررسی رادیوگرافی وجود تومور پانکوست را تأیید می‌کند. معاینه فیزیکی نشان‌دهنده میوز مردمک، پتوز جزئی پلک و آنهیدروز صورت است. کدام یک از ساختارهای زیر به احتمال زیاد آسیب دیده است؟'
$ws.Range("A498").Value = 'هر دو مسیر کمپلمان وابسته به آنتی‌بادی و مستقل از آنتی‌بادی بر کدام جزء کمپلمان همگرا می‌شوند؟'
$ws.Range("A499").Value = 'مثبت بودن شاخص موریس در نوار قلب نشانگر کدام مورد است؟'
$ws.Range("A500").Value = 'بیشترین متابولیسم عبور اول در کدام روش دیده می‌شود:'
$ws.Range("A501").Value = 'شب کوری ممکن است با کدام یک از موارد زیر مرتبط باشد؟ مارس 2004'
$ws.Range("A502").Value = 'سونوگرافی شکمی یک بیمار، کیسه صفرا با ضخیم شدن دیفیوز دیواره و ندول‌های هایپره‌اکو در ناحیه گردن و چند اثر دنباله‌دار ستاره‌ای (comet tail) را نشان داد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A503").Value = 'انتخاب سایه در چه مدت زمانی انجام می‌شود تا از خستگی چشم‌ها جلوگیری شود؟'
$ws.Range("A504").Value = 'بزرگ شدن غدد بناگوشی دو طرفه در همه موارد زیر دیده می‌شود به جز:'
$ws.Range("A505").Value = 'کدام گزینه در مورد کم خونی فقر آهن صحیح است؟  
الف) کم خونی میکروسیتیک هیپوکروم  
ب) کاهش tibc  
ج) افزایش فریتین  
د) کاهش آهن مغز استخوان زودتر از آهن سرم'
$ws.Range("A506").Value = 'در قالبگیری اولیه، سینی استوک باید ________ میلی‌متر بزرگتر از رج باشد.'
$ws.Range("A507").Value = 'پیوند تمام ضخامت در کدام مورد استفاده می‌شود؟'
$ws.Range("A508").Value = 'در حضور دارویی که تمام اثرات نوراپی نفرین و اپی نفرین بر قلب را مسدود می‌کند، سیستم عصبی خودمختار می‌تواند'
$ws.Range("A509").Value = 'شریان سیرکومفلکس ایلیاک سطحی شاخه‌ای از کدام یک از شریان‌های زیر است؟'
$ws.Range("A510").Value = 'در نتیجه آسیب شدید کبدی، متابولیسم اسیدهای آمینه کبدی مختل می‌شود. در این شرایط، سطح پلاسمایی کدام یک از موارد زیر احتمالاً کمتر از حد طبیعی خواهد بود؟'
$ws.Range("A511").Value = 'ksy چیست؟'
$ws.Range("A512").Value = 'بیمار با پاپولووزیکول‌های خارش‌دار مزمن روی سطوح اکستانسور مراجعه می‌کند. کدام یک از موارد زیر نادرست است؟'
$ws.Range("A513").Value = 'هنگام گرفتن قالب پلی سولفاید برای تاج ریختگی، تری های سفارشی معمولاً به دلایل زیر بر تری های آماده ترجیح داده می شوند به جز'
$ws.Range("A514").Value = 'کدام یک از تومورهای سر و گردن زیر بدترین پیش آگهی را دارد؟'
$ws.Range("A515").Value = 'مرکز مهره از چه چیزی تشکیل شده است؟'
$ws.Range("A516").Value = 'قاعده هاچینسون مربوط به کدام مورد است؟'

Write-Output "done"